$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

$ws.Range("A$row").Value = 112170842
$ws.Range("B$row").Value = 96348
$ws.Range("C$row").Value = "Ovaliderad"
$ws.Range("D$row").Value = "VU"
$ws.Range("E$row").Value = 220787
$ws.Range("F$row").Value = "Knärot"
$ws.Range("G$row").Value = "Goodyera repens"
$ws.Range("H$row").Value = "(L.) R. Br."
$ws.Range("I$row").Value = ""
$ws.Range("J$row").Value = ""
$ws.Range("K$row").Value = ""
$ws.Range("L$row").Value = ""
$ws.Range("N$row").Value = ""
$ws.Range("P$row").Value = "Stor Kärmsjön, Ång"
$ws.Range("Q$row").Value = 583983.7355864819
$ws.Range("R$row").Value = 7086628.639155544
$ws.Range("S$row").Value = 25
$ws.Range("T$row").Value = "Västernorrland"
$ws.Range("U$row").Value = "Sollefteå"
$ws.Range("V$row").Value = "Ångermanland"
$ws.Range("W$row").Value = "Junsele"
$ws.Range("Y$row").Value = "'2023-09-15"
$ws.Range("Y$row").Style = "Normal"
$ws.Range("Z$row").Value = "00:00"
$ws.Range("AA$row").Value = "'2023-09-15"
$ws.Range("AA$row").Style = "Normal"
$ws.Range("AB$row").Value = "00:00"
$ws.Range("AD$row").Value = $false
$ws.Range("AE$row").Value = $false
$ws.Range("AF$row").Value = ""
$ws.Range("AG$row").Value = $false
$ws.Range("AT$row").Value = ""
$ws.Range("AW$row").Value = "Maria Johansson"
$ws.Range("AX$row").Value = "Maria Johansson"
$ws.Range("AY$row").Value = ""
